$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-22 Wednesday" "2025-10-23 Thursday"
Replace-Text "629×7=" "402×6="
Replace-Text "427×6=" "914×6="
Replace-Text "678×7=" "410×6="
Replace-Text "679×6=" "199×4="
Replace-Text "394×2=" "940×3="
Replace-Text "339×8=" "549×5="
Replace-Text "272×9=" "275×4="
Replace-Text "919×7=" "461×9="
Replace-Text "685×5=" "665×6="
Replace-Text "729×8=" "853×2="
Replace-Text "142×6=" "330×5="
Replace-Text "288×8=" "580×7="
Replace-Text "746×8=" "294×8="
Replace-Text "300×2=" "408×6="
Replace-Text "331×8=" "484×8="
Replace-Text "443×3=" "846×6="
Replace-Text "770×4=" "995×7="
Replace-Text "612×9=" "695×8="
Replace-Text "890×4=" "461×7="
Replace-Text "796×9=" "342×7="
Replace-Text "778×4=" "587×8="
Replace-Text "396×4=" "705×3="
Replace-Text "624×5=" "395×8="
Replace-Text "999×5=" "987×8="
Replace-Text "537×5=" "214×9="
